$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H3:H9 from "to do" to "done"
$ws.Range("H3:H9").Value = "done"

# Update the active selection to H7 as recorded in the saved view
$ws.Range("H7").Select()
